# repull data, push all data, mean calculation
# Update column F (dSF) values for rows 2-30 (excluding rows 3, 4, 11 which are unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    5  = 1
    6  = 1
    7  = 1
    8  = -4
    9  = 1
    10 = -2
    12 = -2
    13 = -3
    14 = -2
    15 = -1
    16 = -4
    17 = 1
    18 = 1
    19 = -1
    20 = 3
    21 = -4
    22 = 3
    23 = 1
    24 = -2
    25 = -3
    26 = 4
    27 = 1
    28 = -1
    29 = 0
    30 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
